# Adds the 3 new "Electromecanica" game-log rows (41-43) that the
# commit appended to the report sheet, extending the used range from
# A1:Y40 to A1:Y43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ruleta = "C:\Users\jonat\OneDrive\Escritorio\Repositorio\jonatha1992\Predictor_App\Data\Electromecanica.xlsx"
$juego = "No es Simulación"

# Columns, in order A..Y:
# Fecha, Nros jugados, Aciertos Totales, Aciertos Predecidos, V1L, V2L, V3L, V4L,
# l2, dropout, learning rate, epoca, batch_size, Nros a Predecir, Nros Anteriores,
# Cant. Vecinos, Valor_ficha, Limite_juego, Limite_pretendiente, Probabilidad,
# Efectividad, Ruleta, Ganancia, Juego, Predecidos
$newRows = @(
    @("2024-10-16 22:53:42", "", 13, 4, 9, 0, 0, 0, "", "", "", "", "", 10, 10, 1, "", 5, "", 20, "26%", $ruleta, "", $juego, 50),
    @("2024-10-19 00:27:08", "", 2, 0, 2, 0, 0, 0, "", "", "", "", "", 10, 10, 1, "", 5, "", 20, "13%", $ruleta, "", $juego, 15),
    @("2024-10-19 00:30:29", "", 13, 4, 9, 0, 0, 0, "", "", "", "", "", 10, 10, 1, "", 5, "", 20, "34%", $ruleta, "", $juego, 38)
)

$startRow = 41
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $val = $values[$c - 1]
        $cell = $ws.Cells.Item($r, $c)
        if ($val -is [string]) {
            # Force literal text so Excel doesn't reinterpret an empty
            # string as a cleared cell or a "NN%" string as a percent
            # number - mirrors the source file, where these columns are
            # stored as inline strings, not numbers.
            if ($val -eq "" -or $val -match "^\d+(\.\d+)?%$") {
                $cell.Value = "'" + $val
            }
            else {
                $cell.Value = $val
            }
        }
        else {
            $cell.Value = $val
        }
    }
}

Write-Host "Added rows 41-43 (Electromecanica)"
